$d = $word.ActiveDocument

# 1) Owner organization name field (first owner) - noData -> notApplicable
$d.Content.Find.Execute(
    "{d.parcels[i].owners[i].organizationName:ifEM():show(.noData)}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{d.parcels[i].owners[i].organizationName:ifEM():show(.notApplicable)}", 2)

# 2) Owner organization name field (second owner / transferee) - noData -> notApplicable
$d.Content.Find.Execute(
    "{d.parcels[i].owners[i+1].organizationName:ifEM():show(.noData)}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{d.parcels[i].owners[i+1].organizationName:ifEM():show(.notApplicable)}", 2)

# 3) soilFillTypeToPlace merge field - merge the split runs back into one run
$d.Content.Find.Execute(
    "{d.soilFillTypeToPlace:convCRLF:ifEM():show(.noData)}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{d.soilFillTypeToPlace:convCRLF:ifEM():show(.noData)}", 2)

# 4) soilProjectDuration merge field - merge the split runs back into one run
$d.Content.Find.Execute(
    "{d.soilProjectDuration:ifEM():show(.noData)}",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "{d.soilProjectDuration:ifEM():show(.noData)}", 2)

# 5) Footer "Page X of Y" cached PAGE field result: 3 -> 7
#    (scoped to the footer range, whole-word match so only the standalone
#    cached field-result run "3" is touched, not e.g. a "13")
$sec = $d.Sections.First
$footer = $sec.Footers.Item(1)
$footer.Range.Find.Execute("3", $true, $true, $false, $false, $false, $true, 1, $false, "7", 2)
